$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Collapse the last three multi-value (tabbed) rows into single values ---
# These are the last three rows of the original table (rows 34, 35, 36, 1-based).
$rowCount = $t.Rows.Count
$t.Cell($rowCount - 2, 1).Range.Text = "99.99"
$t.Cell($rowCount - 1, 1).Range.Text = "0"
$t.Cell($rowCount, 1).Range.Text = "18"

# --- Update the first three rows ---
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"

# --- Insert ten new rows right after row 3, each holding one new value ---
# Rows.Add(beforeRow) always inserts immediately before "row 4", so pushing
# values in reverse order (last value first) leaves them in forward order.
$values = @("33", "0.00003", "0.00005", "0.00003", "0.00000", "0.00004", "0.00004", "0.00004", "0.00115", "100.0")
for ($i = $values.Count - 1; $i -ge 0; $i--) {
    $beforeRow = $t.Rows.Item(4)
    $t.Rows.Add($beforeRow) | Out-Null
    $t.Cell(4, 1).Range.Text = $values[$i]
}
